$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had 3 data-bearing rows (header + 2 match rows).
# The first match row (old row 2, "Wilstermann vs Oriente Petrolero",
# Bolivia) is removed entirely; the second match row (old row 3,
# "Atl. San Luis vs U.A.N.L.- Tigres", Mexico) shifts up to become row 2.
$ws.Rows(2).Delete()

# After the row shift, the column order of the "HT correct-score" odds
# block (AW:BD) also changed: "Odd_CS_3-3_HT" moved from the end of the
# block (old BC) to the front (new AW), and the remaining columns
# (old AW:BB) each slide one place to the right (new AX:BC). BD
# ("Odd_CS_4-4_HT") is unchanged. Apply this permutation to both the
# header labels (row 1) and the surviving data row (row 2).
$cols = @("AW", "AX", "AY", "AZ", "BA", "BB", "BC")

foreach ($r in 1, 2) {
    $vals = @()
    foreach ($c in $cols) {
        $vals += , ($ws.Range("$c$r").Value2)
    }
    # new AW = old BC ; new AX..BC = old AW..BB
    $ws.Range("AW$r").Value = $vals[6]
    for ($i = 0; $i -lt 6; $i++) {
        $ws.Range("$($cols[$i + 1])$r").Value = $vals[$i]
    }
}
